$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy H1's formatting (bold font + border + centered alignment)
# onto the new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..33: I gets a constant 1, J mirrors column H
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
